$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (A4 = 21)
$ws.Range("B4").Value = 92.31243281328341
$ws.Range("C4").Value = 15.10956701022704
$ws.Range("D4").Value = 69.04212945962217
$ws.Range("E4").Value = 15.84830353015079
$ws.Range("F4").Value = 84.15169646984921
$ws.Range("G4").Value = 15.84830353015079

# Row 5 (A5 = 27)
$ws.Range("B5").Value = 91.36211995779584
$ws.Range("C5").Value = 3.560091977063498
$ws.Range("D5").Value = 83.02406258426124
$ws.Range("E5").Value = 13.41584543867526
$ws.Range("F5").Value = 86.58415456132474
$ws.Range("G5").Value = 13.41584543867526
$ws.Range("H5").Value = 84.16103557019441

# Row 10 (A10 = 47)
$ws.Range("B10").Value = 95.33920971283109

# Row 11 (A11 = 51)
$ws.Range("B11").Value = 95.35073559458699

# Row 12 (A12 = 57)
$ws.Range("B12").Value = 98.48406223524952

# Row 13 (A13 = 61)
$ws.Range("B13").Value = 75.5966977831971

# Row 15 (A15 = 71)
$ws.Range("B15").Value = 98.26540129188768

# Row 16 (A16 = 77)
$ws.Range("B16").Value = 87.55520630885273

# Row 17 (A17 = 81)
$ws.Range("B17").Value = 88.25445330474211

# Row 20 (A20 = Global)
$ws.Range("B20").Value = 88.76230889935026
$ws.Range("C20").Value = 26.64723565892969
$ws.Range("D20").Value = 48.37613624522655
$ws.Range("E20").Value = 24.97662809584376
$ws.Range("F20").Value = 75.02337190415625
$ws.Range("G20").Value = 24.97662809584376
$ws.Range("H20").Value = 66.67211941872959

# Row 21 footnote text (N21) - add 'Deep Sea' and remove ISSCAAP code 46
$newText = @"
NOTE: Percent coverages are performed across FAO major fishing areas to be consistent with Fishstatj. 
Thus, landings from areas such as 'Salmon', 'Tuna', 'Deep Sea', and 'Sharks' are added back into the FAO major fishing area from where they were reported. 
Percent coverage calculations do not include landings from ISSCAAP codes 61, 62, 63, 64, 71, 72, 73, 74, 81, 82, 83, 91, 92, 93, 94, 
except for stocks from these groups which are included in the assessment.
"@

$ws.Range("N21").Value = $newText
